$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# "Effects" section (rows 25-29): remove the "Add/Remove lines to active layer"
# keybinding block (columns D-F, rows 26-29), and rename the anti-aliasing
# toggle description to flag it as dangerous.
$ws.Range("D26").ClearContents()
$ws.Range("E26").ClearContents()
$ws.Range("F26").ClearContents()

$ws.Range("C27").Value = "anti-aliasing (dangerous)"
$ws.Range("F27").ClearContents()

$ws.Range("D28").ClearContents()
$ws.Range("E28").ClearContents()
$ws.Range("F28").ClearContents()

$ws.Range("F29").ClearContents()

# Move the live selection the way the author last left the sheet.
$ws.Range("C28").Select()
